# Applies the "LinuxForHealth" re-brand edit described by the diff:
#  - Metadata sheet: URL, Version, Date, Publisher updated.
#  - Elements sheet: the ele-1/ext-1 constraint text, previously duplicated
#    on both the base "Extension" row and the "Extension.extension" row,
#    is cleared from the base "Extension" row (AI2) - it now only applies
#    to Extension.extension (AI4), which already has it and is unchanged.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-labor-union"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

# Extension.url's "Fixed Value" cell shared the same shared-string entry as
# the Metadata URL in the source workbook, so it tracks the same rename.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-labor-union"
